# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the report
# has moved from "In Translation" to "Ready for handoff", refreshes the
# related timestamps, and widens the Status columns to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-26 06:37:51"
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-26 06:37:47"
$zhcn.Columns.Item(3).ColumnWidth = 16.33

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-26 06:37:51"
$dede.Columns.Item(3).ColumnWidth = 16.33
